$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new LeetCode entry: 0225 | EASY | Implement Stack using Queues ---
$ws.Range("A69").Value = "0225"
$ws.Range("B69").Value = "EASY"
$ws.Range("C69").Value = "Implement Stack using Queues"
$ws.Range("D69").Value = "method1.cpp"
$ws.Range("E69").Value = "Stack｜Design"
$ws.Range("F69").Value = "too easy.."
$ws.Range("G69").Value = "DONE"
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 100
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 100
$ws.Range("L69").Value = 9.1
$ws.Range("M69").Value = 6.67
$ws.Range("N69").Value = 43847
$ws.Range("O69").Value = 0.14097222222222222

# --- Fix the (slightly stale) formatting on rows 67-68 and apply the same
# standard "data row" formatting (number formats/font/fill/border) to the
# brand-new row 69; values are left untouched by a Formats-only paste. ---
$ws.Range("A65:P65").Copy() | Out-Null
$ws.Range("A67:P69").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# View state: scroll the visible window down near the bottom of the data
# and select the freshly-added row, like a user would after data entry.
$excel.ActiveWindow.ScrollRow = 59
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A69:XFD69").Select()
